$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new (blank) column before column N. This shifts the existing
# "Late" (N) and "Outstanding" (P) columns one place to the right
# (N->O, O->P, P->Q), matching the new Variable Instalments layout.
$ws.Columns("N").Insert()

# Move the selection / view to reflect the new active cell and drop the
# previously scrolled-down view (topLeftCell="A7").
[void]$ws.Range("S6").Select()
